# Generate Report for Handback
#
# This script brings the localization-status workbook up to date after a
# handback: it flips the Overview "Status" column from "Ready for handoff"
# to "Handed back: in sync with en-US", records the generated target/handback
# file names + handback timestamps for both locale sheets (zh-cn, de-de),
# adds hyperlinks on the newly-populated "Latest Target File" column, and
# widens a few columns that now need to show longer text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + column widths
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns so the longer text fits. The host's
# ColumnWidth setter quantizes to 1/6-character steps, so 29.1666... is the
# input that lands on the nearest representable width to the target ~30.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Helper data shared by both locale sheets
# ---------------------------------------------------------------------
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a84df01030434cfc98776d0bb6cbb0db8e92611f/e2e/"
$file85 = "85ed5315-7c25-4f7a-bc30-0d83785a1907.md"
$fileD5 = "d53c0648-273d-495b-85cb-894bd2dbb812.md"

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Column widths: Status (C) and the two newly-filled columns (I, J).
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# Row 2 (85ed5315... source file)
$zhcn.Range("J2").Value = "85ed5315-7c25-4f7a-bc30-0d83785a1907.650e3f9faaf2cca878a9e0cc12d0aafbaddb645c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-02 03:11:15"

# Row 3 (d53c0648... source file)
$zhcn.Range("J3").Value = "d53c0648-273d-495b-85cb-894bd2dbb812.821d33fa18e926573ae3d9a59fdcf5a2fc3296f3.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-02 03:11:15"

# Rebuild the hyperlinks collection so the new "Latest Target File" links
# (I2, I3) sit alongside the existing Source File Name links (A2, A3), in
# display order.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $ghBase + $file85, $null, $null, $file85)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $ghBase + $file85, $null, $null, $file85)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $ghBase + $fileD5, $null, $null, $fileD5)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $ghBase + $fileD5, $null, $null, $fileD5)

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Column widths: Status (C) and the two newly-filled columns (I, J).
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

# Row 2 (85ed5315... source file)
$dede.Range("J2").Value = "85ed5315-7c25-4f7a-bc30-0d83785a1907.650e3f9faaf2cca878a9e0cc12d0aafbaddb645c.de-de.xlf"
$dede.Range("K2").Value = "2016-09-02 03:11:22"

# Row 3 (d53c0648... source file)
$dede.Range("J3").Value = "d53c0648-273d-495b-85cb-894bd2dbb812.821d33fa18e926573ae3d9a59fdcf5a2fc3296f3.de-de.xlf"
$dede.Range("K3").Value = "2016-09-02 03:11:22"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $ghBase + $file85, $null, $null, $file85)
$dede.Hyperlinks.Add($dede.Range("I2"), $ghBase + $file85, $null, $null, $file85)
$dede.Hyperlinks.Add($dede.Range("A3"), $ghBase + $fileD5, $null, $null, $fileD5)
$dede.Hyperlinks.Add($dede.Range("I3"), $ghBase + $fileD5, $null, $null, $fileD5)
